# Update the cached text of the "datetimeFigureOut" date placeholders
# (slide master, notes master, and every slide layout) from 4/27/2011
# to 5/10/2011, mirroring the source document's regenerated date field.

$p = $ppt.ActivePresentation
$oldDate = "4/27/2011"
$newDate = "5/10/2011"

function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster

# Every slide layout belonging to the slide master
$sm = $p.SlideMaster
for ($j = 1; $j -le $sm.CustomLayouts.Count; $j++) {
    Update-DatePlaceholder $sm.CustomLayouts.Item($j)
}

# Notes master date placeholder is only reachable/settable through the
# HeadersFooters object in this object model.
$nm = $p.NotesMaster
if ($nm.HeadersFooters.DateAndTime.Text -eq $oldDate) {
    $nm.HeadersFooters.DateAndTime.Text = $newDate
}
